$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.151.02"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.900.95"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0758"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.74%  "
$ws.Range("D14").Value = "2.177.75"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.733"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("E16").Value = "  +3.76%  "
$ws.Range("D17").Value = "1.919.40"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "35.141.99"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "0.0₃0840"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  +4.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D31").Value = "4.128.72"
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +11.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.51%  "
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("E38").Value = "  -6.97%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("E42").Value = "  +3.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0669"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").Value = "1.305.25"
$ws.Range("E46").Value = "  -3.30%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("E51").Value = "  +6.77%  "
